# Edit summary (see commit message "understand how stage works"):
#   - The last two paragraphs of the diary change:
#       * the old 6th paragraph ("英语三级考试结束。") loses the trailing
#         "_GoBack" bookmark that used to sit right after its text.
#       * a brand-new paragraph "2022年9月20号星期二" is inserted after it,
#         formatted the same way as the other date lines (rFonts hint=eastAsia).
#       * the old trailing empty paragraph gets real text
#         "今天天气好好，适合出去走走。" typed into it (its own paragraph
#         formatting, rFonts hint=default, is left untouched) and the
#         "_GoBack" bookmark is now placed at the end of that new text.
#
# Because the sandboxed engine does not automatically stamp freshly-typed
# runs with the "eastAsia" rFonts hint the way real Word's IME does, new
# CJK runs are produced here by duplicating the *formatting* of an existing,
# already-correctly-hinted run (paragraph 1's) via Range.FormattedText, and
# then swapping in the desired text with a scoped Find/Replace (which keeps
# the run's rPr intact).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark (currently collapsed right
#    after "英语三级考试结束。", i.e. at the end of paragraph 6).
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# 2. Insert a brand-new paragraph, formatted like paragraph 1 (rFonts
#    hint="eastAsia" on both the paragraph mark and the run), right before
#    the trailing empty paragraph. Start from a full copy of paragraph 1
#    (text + mark) so the new paragraph mark/pPr get the right formatting,
#    then overwrite the text.
# ---------------------------------------------------------------------
$templatePara = $d.Paragraphs.Item(1)
$trailingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPos = $trailingPara.Range.Start
$newParaSpot = $d.Range($insertPos, $insertPos)
$newParaSpot.FormattedText = $templatePara.Range.FormattedText

$datePara = $d.Paragraphs.Item(7)
$datePara.Range.Find.Execute("2022年9月16日星期五", $false, $false, $false, $false, $false, $true, 1, $false, "2022年9月20号星期二", 2)

# ---------------------------------------------------------------------
# 3. Give the (now last) originally-empty paragraph a real CJK run by
#    duplicating just the formatted run (no paragraph mark) from
#    paragraph 1, then overwrite its text. The paragraph's own pPr
#    (rFonts hint="default") is left exactly as it was.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$templatePara2 = $d.Paragraphs.Item(1)
$templateRun = $d.Range($templatePara2.Range.Start, $templatePara2.Range.End - 1)
$runSpot = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$runSpot.FormattedText = $templateRun.FormattedText

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Find.Execute("2022年9月16日星期五", $false, $false, $false, $false, $false, $true, 1, $false, "今天天气好好，适合出去走走。", 2)

# ---------------------------------------------------------------------
# 4. Re-create the "_GoBack" bookmark at the end of the new last
#    paragraph (collapsed, right before its paragraph mark).
#
#    Quirk: adding a collapsed bookmark whose position lands exactly on a
#    paragraph-end boundary gets mis-resolved by this environment. Work
#    around it by anchoring the bookmark while a placeholder run still
#    separates it from the paragraph mark, then deleting the placeholder
#    (ordinary bookmark range-adjustment keeps the bookmark glued to the
#    text that now precedes it).
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$markPos = $lastPara.Range.End - 1
$placeholder = $d.Range($markPos, $markPos)
$placeholder.InsertBefore("PLACEHOLDER")

$bmSpot = $d.Range($markPos, $markPos)
$d.Bookmarks.Add("_GoBack", $bmSpot)

$goBack2 = $d.Bookmarks.Item("_GoBack")
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$placeholderRange = $d.Range($goBack2.End, $lastPara.Range.End - 1)
$placeholderRange.Delete()
